$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric (e.g. "28.422.44") - force text format so
# Excel does not coerce them into numbers, matching the source inline strings.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "28.422.44"
$ws.Cells.Item(2, 5).Value = "  +0.18%  "
$ws.Cells.Item(3, 4).Value = "1.819.64"
$ws.Cells.Item(3, 5).Value = "  -0.41%  "
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$ws.Cells.Item(5, 4).Value = "315.22"
$ws.Cells.Item(5, 5).Value = "  -0.71%  "
$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  -0.01%  "
$ws.Cells.Item(7, 4).Value = "0.5128"
$ws.Cells.Item(7, 5).Value = "  -4.20%  "
$ws.Cells.Item(8, 4).Value = "0.3917"
$ws.Cells.Item(8, 5).Value = "  -3.66%  "
$ws.Cells.Item(9, 4).Value = "0.07873"
$ws.Cells.Item(9, 5).Value = "  +3.63%  "
$ws.Cells.Item(10, 4).Value = "41.75"
$ws.Cells.Item(10, 5).Value = "  -0.19%  "
$ws.Cells.Item(11, 4).Value = "1.107"
$ws.Cells.Item(11, 5).Value = "  +0.30%  "
$ws.Cells.Item(12, 4).Value = "20.98"
$ws.Cells.Item(12, 5).Value = "  +1.13%  "
$ws.Cells.Item(13, 4).Value = "6.253"
$ws.Cells.Item(13, 5).Value = "  -1.19%  "
$ws.Cells.Item(14, 4).Value = "1.001"
$ws.Cells.Item(14, 5).Value = "  +0.02%  "
$ws.Cells.Item(15, 4).Value = "7.489"
$ws.Cells.Item(15, 5).Value = "  -1.18%  "
$ws.Cells.Item(16, 4).Value = "1.817.30"
$ws.Cells.Item(16, 5).Value = "  -0.55%  "
$ws.Cells.Item(17, 4).Value = "0.00001125"
$ws.Cells.Item(17, 5).Value = "  +4.95%  "
$ws.Cells.Item(18, 4).Value = "92.61"
$ws.Cells.Item(18, 5).Value = "  +3.72%  "
$ws.Cells.Item(19, 4).Value = "0.06618"
$ws.Cells.Item(19, 5).Value = "  +0.10%  "
$ws.Cells.Item(20, 4).Value = "17.69"
$ws.Cells.Item(20, 5).Value = "  +0.49%  "
$ws.Cells.Item(21, 4).Value = "1.001"
$ws.Cells.Item(21, 5).Value = "  -0.01%  "
$ws.Cells.Item(22, 4).Value = "6.089"
$ws.Cells.Item(22, 5).Value = "  +0.10%  "
$ws.Cells.Item(23, 4).Value = "28.449.26"
$ws.Cells.Item(23, 5).Value = "  +0.19%  "
$ws.Cells.Item(24, 4).Value = "11.23"
$ws.Cells.Item(24, 5).Value = "  +0.46%  "
$ws.Cells.Item(25, 4).Value = "2.270"
$ws.Cells.Item(25, 5).Value = "  +3.99%  "
$ws.Cells.Item(26, 4).Value = "21.03"
$ws.Cells.Item(26, 5).Value = "  +2.23%  "
$ws.Cells.Item(27, 4).Value = "2.029.21"
$ws.Cells.Item(27, 5).Value = "  -0.48%  "
$ws.Cells.Item(28, 4).Value = "154.96"
$ws.Cells.Item(28, 5).Value = "  -1.86%  "
$ws.Cells.Item(29, 4).Value = "2.399"
$ws.Cells.Item(29, 5).Value = "  -2.74%  "
$ws.Cells.Item(30, 4).Value = "125.64"
$ws.Cells.Item(30, 5).Value = "  +1.63%  "
$ws.Cells.Item(31, 4).Value = "0.1098"
$ws.Cells.Item(31, 5).Value = "  +0.54%  "
$ws.Cells.Item(32, 4).Value = "1.104"
$ws.Cells.Item(32, 5).Value = "  -1.65%  "
$ws.Cells.Item(33, 4).Value = "5.673"
$ws.Cells.Item(33, 5).Value = "  +0.45%  "
$ws.Cells.Item(34, 4).Value = "3.650"
$ws.Cells.Item(34, 5).Value = "  +0.23%  "
$ws.Cells.Item(35, 4).Value = "0.07059"
$ws.Cells.Item(35, 5).Value = "  -2.52%  "
$ws.Cells.Item(36, 4).Value = "0.2215"
$ws.Cells.Item(36, 5).Value = "  -1.49%  "
$ws.Cells.Item(37, 4).Value = "0.02324"
$ws.Cells.Item(37, 5).Value = "  -0.45%  "
$ws.Cells.Item(38, 4).Value = "5.183"
$ws.Cells.Item(38, 5).Value = "  -0.29%  "
$ws.Cells.Item(39, 4).Value = "8.781"
$ws.Cells.Item(39, 5).Value = "  -0.59%  "
$ws.Cells.Item(40, 4).Value = "0.6249"
$ws.Cells.Item(40, 5).Value = "  -0.24%  "
$ws.Cells.Item(41, 4).Value = "11.26"
$ws.Cells.Item(41, 5).Value = "  -0.18%  "
$ws.Cells.Item(42, 4).Value = "1.176"
$ws.Cells.Item(42, 5).Value = "  -0.54%  "
$ws.Cells.Item(43, 4).Value = "1.0000"
$ws.Cells.Item(43, 5).Value = "  -0.08%  "
$ws.Cells.Item(44, 4).Value = "1.396"
$ws.Cells.Item(44, 5).Value = "  -0.40%  "
$ws.Cells.Item(45, 4).Value = "13.46"
$ws.Cells.Item(45, 5).Value = "  -0.15%  "
$ws.Cells.Item(46, 4).Value = "3.736"
$ws.Cells.Item(46, 5).Value = "  +0.84%  "
$ws.Cells.Item(47, 4).Value = "0.5892"
$ws.Cells.Item(47, 5).Value = "  +0.79%  "
$ws.Cells.Item(48, 4).Value = "124.59"
$ws.Cells.Item(48, 5).Value = "  -0.68%  "
$ws.Cells.Item(49, 4).Value = "1.969"
$ws.Cells.Item(49, 5).Value = "  -0.73%  "
$ws.Cells.Item(50, 4).Value = "1.190"
$ws.Cells.Item(50, 5).Value = "  -1.08%  "
$ws.Cells.Item(51, 4).Value = "0.06888"
$ws.Cells.Item(51, 5).Value = "  +0.02%  "
